$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.588.69"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "3.396.30"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").Value = "560.11"
$ws.Range("E5").Value = "  -0.15%  "

$ws.Range("D6").Value = "175.89"
$ws.Range("E6").Value = "  -0.43%  "

$ws.Range("D7").Value = "0.629"
$ws.Range("E7").Value = "  +1.15%  "

$ws.Range("D8").Value = "3.385.70"
$ws.Range("E8").Value = "  +0.84%  "

$ws.Range("E9").Value = "  +0.03%  "

$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +4.10%  "

$ws.Range("D11").Value = "0.637"
$ws.Range("E11").Value = "  +0.87%  "

$ws.Range("D12").Value = "53.84"
$ws.Range("E12").Value = "  -2.30%  "

$ws.Range("D13").Value = "0.0000280"
$ws.Range("E13").Value = "  +1.17%  "

$ws.Range("D14").Value = "9.24"
$ws.Range("E14").Value = "  +1.34%  "

$ws.Range("D15").Value = "3.937.34"
$ws.Range("E15").Value = "  +0.85%  "

$ws.Range("D16").Value = "18.36"
$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  +1.76%  "

$ws.Range("D18").Value = "3.383.08"
$ws.Range("E18").Value = "  +0.52%  "

$ws.Range("D19").Value = "65.430.77"
$ws.Range("E19").Value = "  +1.34%  "

$ws.Range("D20").Value = "11.90"
$ws.Range("E20").Value = "  +0.23%  "

$ws.Range("E21").Value = "  +1.32%  "

$ws.Range("D22").Value = "464.47"
$ws.Range("E22").Value = "  +0.54%  "

$ws.Range("D23").Value = "4.99"
$ws.Range("E23").Value = "  +3.39%  "

$ws.Range("D24").Value = "4.14"
$ws.Range("E24").Value = "  +0.66%  "

$ws.Range("D25").Value = "14.36"
$ws.Range("E25").Value = "  +6.23%  "

$ws.Range("D26").Value = "87.49"
$ws.Range("E26").Value = "  +1.50%  "

$ws.Range("E27").Value = "  +3.05%  "

$ws.Range("D28").Value = "10.75"
$ws.Range("E28").Value = "  -1.24%  "

$ws.Range("D29").Value = "8.75"
$ws.Range("E29").Value = "  -0.71%  "

$ws.Range("D30").Value = "31.13"
$ws.Range("E30").Value = "  +3.23%  "

$ws.Range("D31").Value = "6.56"
$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("D32").Value = "63.63"
$ws.Range("E32").Value = "  +6.71%  "

$ws.Range("D33").Value = "11.53"
$ws.Range("E33").Value = "  +0.21%  "

$ws.Range("D34").Value = "578.50"
$ws.Range("E34").Value = "  -0.29%  "

$ws.Range("E35").Value = "  -0.22%  "

$ws.Range("E36").Value = "  +0.00%  "

$ws.Range("D37").Value = "3.60"
$ws.Range("E37").Value = "  +4.12%  "

$ws.Range("D38").Value = "0.142"
$ws.Range("E38").Value = "  +0.99%  "

$ws.Range("D39").Value = "36.01"
$ws.Range("E39").Value = "  +0.16%  "

$ws.Range("D40").Value = "0.375"
$ws.Range("E40").Value = "  +0.45%  "

$ws.Range("D41").Value = "0.0₃0744"
$ws.Range("E41").Value = "  -2.10%  "

$ws.Range("D42").Value = "3.113.78"
$ws.Range("E42").Value = "  +0.49%  "

$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").Value = "2.80"
$ws.Range("E43").Value = "  -1.49%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "0.0418"
$ws.Range("E44").Value = "  +1.23%  "

$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.134"
$ws.Range("E45").Value = "  +1.82%  "

$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").Value = "2.45"
$ws.Range("E46").Value = "  -2.89%  "

$ws.Range("D47").Value = "3.16"
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  +0.00%  "

$ws.Range("D49").Value = "140.45"
$ws.Range("E49").Value = "  +3.35%  "

$ws.Range("D50").Value = "2.55"
$ws.Range("E50").Value = "  -1.69%  "

$ws.Range("D51").Value = "8.47"
$ws.Range("E51").Value = "  +1.16%  "
